$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.714.49'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = '1.600.83'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.54'
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('E6').Value = '  -0.57%  '
$ws.Range('E7').Value = '  +0.32%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.246'
$ws.Range('E9').Value = '  -0.16%  '
$ws.Range('E10').Value = '  +0.84%  '
$ws.Range('E11').Value = '  +0.42%  '
$ws.Range('D12').Value = '1.826.52'
$ws.Range('E12').Value = '  +0.20%  '
$ws.Range('D13').Value = '1.610.85'
$ws.Range('E13').Value = '  +0.93%  '
$ws.Range('E14').Value = '  +0.28%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.523'
$ws.Range('E15').Value = '  +0.12%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.86'
$ws.Range('E16').Value = '  +0.12%  '
$ws.Range('D17').Value = '26.688.47'
$ws.Range('E17').Value = '  +0.27%  '
$ws.Range('D18').Value = '0.0₃0743'
$ws.Range('E18').Value = '  +0.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '210.12'
$ws.Range('E19').Value = '  +0.89%  '
$ws.Range('E20').Value = '  +0.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.16'
$ws.Range('E21').Value = '  +2.27%  '
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.28'
$ws.Range('E23').Value = '  -2.39%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.94'
$ws.Range('E24').Value = '  +0.41%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.23'
$ws.Range('E25').Value = '  -0.84%  '
$ws.Range('E26').Value = '  +0.34%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.09'
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('E28').Value = '  -0.78%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.36'
$ws.Range('E29').Value = '  +0.34%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0511'
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('E32').Value = '  +0.98%  '
$ws.Range('E33').Value = '  +0.80%  '
$ws.Range('D34').Value = '1.297.43'
$ws.Range('E34').Value = '  +1.53%  '
$ws.Range('E35').Value = '  +0.72%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.604'
$ws.Range('E36').Value = '  -2.44%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.49'
$ws.Range('E37').Value = '  +0.36%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.17'
$ws.Range('E38').Value = '  +8.87%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0170'
$ws.Range('E39').Value = '  -0.69%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.828'
$ws.Range('E40').Value = '  -1.24%  '
$ws.Range('E41').Value = '  -1.48%  '
$ws.Range('E42').Value = '  -0.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.778'
$ws.Range('E43').Value = '  -0.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.99'
$ws.Range('E44').Value = '  -1.69%  '
$ws.Range('D45').Value = '1.739.15'
$ws.Range('E45').Value = '  +0.09%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '90.53'
$ws.Range('E47').Value = '  -2.88%  '
$ws.Range('E48').Value = '  -0.74%  '
$ws.Range('E49').Value = '  +1.73%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.42'
$ws.Range('E51').Value = '  -0.20%  '
